$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "Jan 2023 - Dec 2023 (16/04/24)"
$ws.Range("D2").Value = "Apr 2023 - Mar 2024 (18/07/24)"
$ws.Range("C3").Value = "Jan 2023 - Dec 2023 (16/04/24)"
$ws.Range("D3").Value = "Apr 2023 - Mar 2024 (18/07/24)"
$ws.Range("C4").Value = "Jan 2023 - Dec 2023 (16/04/24)"
$ws.Range("D4").Value = "Apr 2023 - Mar 2024 (18/07/24)"

$ws.Range("D4").Select()
